$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.777.89'
$ws.Range("E2").Value = '  -1.20%  '
$ws.Range("D3").Value = '1.598.48'
$ws.Range("E3").Value = '  -1.97%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '210.62'
$ws.Range("E5").Value = '  -2.47%  '
$ws.Range("E6").Value = '  -1.81%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -2.01%  '
$ws.Range("D9").Value = '0.0616'
$ws.Range("E9").Value = '  -1.03%  '
$ws.Range("D10").Value = '19.61'
$ws.Range("E10").Value = '  -2.51%  '
$ws.Range("D11").Value = '0.0839'
$ws.Range("E11").Value = '  -1.13%  '
$ws.Range("D12").Value = '1.807.70'
$ws.Range("E12").Value = '  -2.74%  '
$ws.Range("D13").Value = '1.613.30'
$ws.Range("E13").Value = '  -1.15%  '
$ws.Range("E14").Value = '  -1.45%  '
$ws.Range("D15").Value = '0.528'
$ws.Range("E15").Value = '  -2.16%  '
$ws.Range("D16").Value = '26.731.24'
$ws.Range("E16").Value = '  -1.27%  '
$ws.Range("D17").Value = '63.41'
$ws.Range("E17").Value = '  -3.41%  '
$ws.Range("D18").Value = '0.0₃0727'
$ws.Range("E18").Value = '  -0.65%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '208.70'
$ws.Range("E19").Value = '  -2.38%  '
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = '1.00'
$ws.Range("E20").Value = '  +0.11%  '
$ws.Range("D21").Value = '6.73'
$ws.Range("E21").Value = '  -1.19%  '
$ws.Range("D22").Value = '4.27'
$ws.Range("E22").Value = '  -2.46%  '
$ws.Range("D23").Value = '2.33'
$ws.Range("E23").Value = '  -6.78%  '
$ws.Range("E24").Value = '  -2.83%  '
$ws.Range("D25").Value = '146.17'
$ws.Range("E25").Value = '  -0.60%  '
$ws.Range("D26").Value = '7.48'
$ws.Range("E26").Value = '  +1.35%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("E28").Value = '  -4.93%  '
$ws.Range("D29").Value = '15.27'
$ws.Range("E29").Value = '  -1.99%  '
$ws.Range("D30").Value = '0.0499'
$ws.Range("E30").Value = '  -0.90%  '
$ws.Range("E31").Value = '  -2.70%  '
$ws.Range("D32").Value = '3.25'
$ws.Range("E32").Value = '  -3.06%  '
$ws.Range("E33").Value = '  +23.38%  '
$ws.Range("E34").Value = '  -2.14%  '
$ws.Range("D35").Value = '1.312.49'
$ws.Range("E35").Value = '  +0.75%  '
$ws.Range("D36").Value = '1.52'
$ws.Range("E36").Value = '  -2.95%  '
$ws.Range("D37").Value = '2.43'
$ws.Range("E37").Value = '  -0.65%  '
$ws.Range("D38").Value = '0.0173'
$ws.Range("E38").Value = '  -1.47%  '
$ws.Range("D39").Value = '0.820'
$ws.Range("E39").Value = '  -2.50%  '
$ws.Range("E40").Value = '  +0.10%  '
$ws.Range("D41").Value = '0.787'
$ws.Range("E41").Value = '  -2.24%  '
$ws.Range("E42").Value = '  -3.61%  '
$ws.Range("D43").Value = '5.26'
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("D44").Value = '62.75'
$ws.Range("E44").Value = '  +0.89%  '
$ws.Range("D45").Value = '1.736.42'
$ws.Range("E45").Value = '  -1.68%  '
$ws.Range("D46").Value = '88.78'
$ws.Range("E46").Value = '  -1.96%  '
$ws.Range("D47").Value = '1.60'
$ws.Range("E47").Value = '  +0.47%  '
$ws.Range("D48").Value = '0.817'
$ws.Range("E48").Value = '  +1.12%  '
$ws.Range("D49").Value = '0.0509'
$ws.Range("E49").Value = '  -0.87%  '
$ws.Range("D50").Value = '0.0974'
$ws.Range("E50").Value = '  +2.84%  '
$ws.Range("D51").Value = '0.0₇0956'
$ws.Range("E51").Value = '  -10.51%  '
